$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F (Maintainer) ---
$ws.Columns.Item(6).ColumnWidth = 15.5

# Header
$ws.Range("F1").Value = "Maintainer"

# Row 10 - ggbiplot maintainer (hyperlinked email)
$ws.Range("F10").Value = "friendly@yorku.ca"
$ws.Hyperlinks.Add($ws.Range("F10"), "mailto:friendly@yorku.ca") | Out-Null
$ws.Range("E10").Copy()
$ws.Range("F10").PasteSpecial(-4122)

# Row 4 - adegraphics maintainer
$ws.Range("F4").Value = "stephane.dray@univ-lyon1.fr"

# Row 6 - biplotEZ maintainer (hyperlinked email)
$ws.Range("F6").Value = "muvisu@sun.ac.za"
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:muvisu@sun.ac.za") | Out-Null
$ws.Range("E6").Copy()
$ws.Range("F6").PasteSpecial(-4122)

# Row 8 - FactoMineR / factoextra maintainer (hyperlinked email)
$ws.Range("F8").Value = "francois.husson@institut-agro.fr"
$ws.Hyperlinks.Add($ws.Range("F8"), "mailto:francois.husson@institut-agro.fr") | Out-Null
$ws.Range("E8").Copy()
$ws.Range("F8").PasteSpecial(-4122)

# Row 9 - factoextra maintainer
$ws.Range("F9").Value = "alboukadel.kassambara@gmail.com"

# Row 13 - new "ordr" package row (value first, hyperlink added later to match rId order)
$ws.Range("A13").Value = "ordr"

# Row 11 - MultBiplotR maintainer
$ws.Range("F11").Value = "villardon@usal.es"

# Row 12 - pcaMethods maintainer (hyperlinked email)
$ws.Range("F12").Value = "henning.red@gmail.com"
$ws.Hyperlinks.Add($ws.Range("F12"), "mailto:henning.red@gmail.com") | Out-Null
$ws.Range("E12").Copy()
$ws.Range("F12").PasteSpecial(-4122)

# Row 13 continued - maintainer email for "ordr"
$ws.Range("F13").Value = "cornelioid@gmail.com"

# Row 13 - vignette link (value first, hyperlink added after A13 per rId order)
$ws.Range("E13").Value = "Ordination in the tidyverse"

# Row 13 - "biplot" column re-using existing "ggplot" string
$ws.Range("D13").Value = "ggplot"

# Now register the remaining hyperlinks in the exact order used by the workbook
$ws.Hyperlinks.Add($ws.Range("A13"), "https://cran.r-project.org/package=ordr") | Out-Null
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("E13"), "https://cornelioid.github.io/ordr/articles/ordr.html") | Out-Null
$ws.Range("E12").Copy()
$ws.Range("E13").PasteSpecial(-4122)

# Update selection to reflect final cursor position
$ws.Range("F15").Select()
